$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 201.625
$ws.Range("I11").Value = 201.625
$ws.Range("K11").Value = 201.625
$ws.Range("M11").Value = -61.625
$ws.Range("H137").Value = 987.1
$ws.Range("I137").Value = 964.38464
$ws.Range("J137").Value = 1001.6229
$ws.Range("K137").Value = 2893.15392
$ws.Range("L137").Value = 3004.8687
$ws.Range("M137").Value = -343.1539199999997
$ws.Range("N137").Value = -8104.8687
$ws.Range("H138").Value = 5559350
$ws.Range("I138").Value = 2605.625
$ws.Range("J138").Value = 10004746
$ws.Range("K138").Value = 7816.875
$ws.Range("L138").Value = 30014238
$ws.Range("M138").Value = -2676.875
$ws.Range("N138").Value = -30024518

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 93.333336
$ws.Range("I5").Value = 95
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 95
$ws.Range("L5").Value = 90
$ws.Range("M5").Value = 17
$ws.Range("N5").Value = -314
$ws.Range("H32").Value = 4254.933
$ws.Range("I32").Value = 3553.838
$ws.Range("J32").Value = 7497.5
$ws.Range("K32").Value = 3553.838
$ws.Range("L32").Value = 7497.5
$ws.Range("M32").Value = -3266.838
$ws.Range("N32").Value = -8071.5
$ws.Range("H43").Value = 9246.666999999999
$ws.Range("I43").Value = 9342
$ws.Range("J43").Value = 9199
$ws.Range("K43").Value = 9342
$ws.Range("L43").Value = 9199
$ws.Range("M43").Value = -9029
$ws.Range("N43").Value = -9825
$ws.Range("H63").Value = 3461.5386
$ws.Range("I63").Value = 2457.1428
$ws.Range("J63").Value = 4633.3335
$ws.Range("K63").Value = 2457.1428
$ws.Range("L63").Value = 4633.3335
$ws.Range("M63").Value = -1771.1428
$ws.Range("N63").Value = -6005.3335
$ws.Range("H66").Value = 3461.5386
$ws.Range("I66").Value = 2457.1428
$ws.Range("J66").Value = 4633.3335
$ws.Range("K66").Value = 12285.714
$ws.Range("L66").Value = 23166.6675
$ws.Range("M66").Value = -8853.714
$ws.Range("N66").Value = -30030.6675
$ws.Range("H97").Value = 1029.4117
$ws.Range("I97").Value = 929.8333
$ws.Range("J97").Value = 1083.7273
$ws.Range("K97").Value = 929.8333
$ws.Range("L97").Value = 1083.7273
$ws.Range("M97").Value = -433.8333
$ws.Range("N97").Value = -2075.7273
$ws.Range("H122").Value = 2128.5715
$ws.Range("I122").Value = 2150
$ws.Range("J122").Value = 2120
$ws.Range("K122").Value = 6450
$ws.Range("L122").Value = 6360
$ws.Range("M122").Value = -4000
$ws.Range("N122").Value = -11260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 93.333336
$ws.Range("I4").Value = 95
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 95
$ws.Range("L4").Value = 90
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = -320
$ws.Range("H5").Value = 13000
$ws.Range("I5").Value = 6000
$ws.Range("J5").Value = 20000
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = -5887
$ws.Range("N5").Value = -20226
$ws.Range("H62").Value = 28000
$ws.Range("J62").Value = 28000
$ws.Range("L62").Value = 28000
$ws.Range("N62").Value = -29372
$ws.Range("H65").Value = 28000
$ws.Range("J65").Value = 28000
$ws.Range("L65").Value = 84000
$ws.Range("N65").Value = -90864
$ws.Range("H86").Value = 2216.5833
$ws.Range("I86").Value = 2145.3635
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 2145.3635
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1022.3635
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 2216.5833
$ws.Range("I89").Value = 2145.3635
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 10726.8175
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -5110.817499999999
$ws.Range("N89").Value = -26232
$ws.Range("H99").Value = 2411.8
$ws.Range("I99").Value = 2482.5334
$ws.Range("J99").Value = 2199.6
$ws.Range("K99").Value = 2482.5334
$ws.Range("L99").Value = 2199.6
$ws.Range("M99").Value = -984.5333999999998
$ws.Range("N99").Value = -5195.6
$ws.Range("H105").Value = 2448.2415
$ws.Range("I105").Value = 2406.1875
$ws.Range("K105").Value = 2406.1875
$ws.Range("M105").Value = -659.1875
$ws.Range("H107").Value = 2533.4
$ws.Range("I107").Value = 2571.6428
$ws.Range("J107").Value = 1998
$ws.Range("K107").Value = 2571.6428
$ws.Range("L107").Value = 1998
$ws.Range("M107").Value = -651.6428000000001
$ws.Range("N107").Value = -5838

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 138.05556
$ws.Range("I7").Value = 129.28572
$ws.Range("J7").Value = 143.63637
$ws.Range("K7").Value = 129.28572
$ws.Range("L7").Value = 143.63637
$ws.Range("M7").Value = -16.28572
$ws.Range("N7").Value = -369.63637
$ws.Range("H8").Value = 4880
$ws.Range("J8").Value = 5950
$ws.Range("L8").Value = 5950
$ws.Range("N8").Value = -6230
$ws.Range("H22").Value = 531.4286
$ws.Range("I22").Value = 389.2
$ws.Range("J22").Value = 887
$ws.Range("K22").Value = 389.2
$ws.Range("L22").Value = 887
$ws.Range("M22").Value = -39.19999999999999
$ws.Range("N22").Value = -1587
$ws.Range("H47").Value = 27714
$ws.Range("J47").Value = 27714
$ws.Range("L47").Value = 27714
$ws.Range("N47").Value = -28846
$ws.Range("H140").Value = 65475.715
$ws.Range("J140").Value = 65475.715
$ws.Range("L140").Value = 65475.715
$ws.Range("N140").Value = -75835.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 1969.2778
$ws.Range("J21").Value = 1668.421
$ws.Range("L21").Value = 5005.263
$ws.Range("N21").Value = -5351.263
$ws.Range("H48").Value = 1951.3043
$ws.Range("J48").Value = 1951.3043
$ws.Range("L48").Value = 5853.9129
$ws.Range("N48").Value = -6353.9129
$ws.Range("H131").Value = 863.9167
$ws.Range("I131").Value = 519.4286
$ws.Range("J131").Value = 891.0112
$ws.Range("K131").Value = 1558.2858
$ws.Range("L131").Value = 2673.0336
$ws.Range("M131").Value = 3481.7142
$ws.Range("N131").Value = -12753.0336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 17575
$ws.Range("J63").Value = 17575
$ws.Range("L63").Value = 17575
$ws.Range("N63").Value = -18947
$ws.Range("H66").Value = 17575
$ws.Range("J66").Value = 17575
$ws.Range("L66").Value = 52725
$ws.Range("N66").Value = -59589
$ws.Range("H82").Value = 28699.666
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 28699.666
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 28699.666
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -29465.666
$ws.Range("H85").Value = 28699.666
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 28699.666
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 28699.666
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -31351.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1513.9
$ws.Range("I16").Value = 1367.6154
$ws.Range("J16").Value = 1785.5714
$ws.Range("K16").Value = 1367.6154
$ws.Range("L16").Value = 1785.5714
$ws.Range("M16").Value = -1197.6154
$ws.Range("N16").Value = -2125.5714
$ws.Range("H21").Value = 56671.332
$ws.Range("J21").Value = 56671.332
$ws.Range("L21").Value = 56671.332
$ws.Range("N21").Value = -57019.332
$ws.Range("H30").Value = 15603.6
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 19254.5
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 19254.5
$ws.Range("M30").Value = -892
$ws.Range("N30").Value = -19470.5
$ws.Range("H35").Value = 5415.5
$ws.Range("I35").Value = 5415.5
$ws.Range("K35").Value = 5415.5
$ws.Range("M35").Value = -5079.5
$ws.Range("H80").Value = 24000
$ws.Range("J80").Value = 24000
$ws.Range("L80").Value = 24000
$ws.Range("N80").Value = -26246
$ws.Range("H83").Value = 24000
$ws.Range("J83").Value = 24000
$ws.Range("L83").Value = 72000
$ws.Range("N83").Value = -83232

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 8998.1
$ws.Range("I55").Value = 5349.3335
$ws.Range("J55").Value = 10561.857
$ws.Range("K55").Value = 5349.3335
$ws.Range("L55").Value = 10561.857
$ws.Range("M55").Value = -5072.3335
$ws.Range("N55").Value = -11115.857
